$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 39
$ws.Range("D2").Value = 3.81986106537417
$ws.Range("E2").Value = 2.586039266456233
$ws.Range("F2").Value = 0.4140976934050983
$ws.Range("G2").Value = 0.8382969538186891
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 0.447415760511841
$ws.Range("E3").Value = 0.2150763429713992
$ws.Range("F3").Value = 0.1075381714856996
$ws.Range("G3").Value = 0.3422344565216304
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = 1762.94051932889
$ws.Range("E4").Value = 1129.625876877023
$ws.Range("F4").Value = 180.8849061547703
$ws.Range("G4").Value = 366.1823483594814
$ws.Range("C5").Value = 39
$ws.Range("D5").Value = 32.86324786324786
$ws.Range("E5").Value = 32.72235835850176
$ws.Range("F5").Value = 5.239770832095349
$ws.Range("G5").Value = 10.60736149272999
$ws.Range("C6").Value = 39
$ws.Range("D6").Value = 15.52991452991453
$ws.Range("E6").Value = 8.081716405008416
$ws.Range("F6").Value = 1.294110327510283
$ws.Range("G6").Value = 2.619789394470025
$ws.Range("C7").Value = 39
$ws.Range("D7").Value = 2.166337988024551
$ws.Range("E7").Value = 6.582499674035399
$ws.Range("F7").Value = 1.054043520225876
$ws.Range("G7").Value = 2.133799550854492
$ws.Range("C8").Value = 39
$ws.Range("D8").Value = 60.51282051282051
$ws.Range("E8").Value = 34.77612999276881
$ws.Range("F8").Value = 5.568637492227787
$ws.Range("G8").Value = 11.27311724020732
$ws.Range("C9").Value = 39
$ws.Range("D9").Value = 0.6410256410256411
$ws.Range("E9").Value = 4.003203845127178
$ws.Range("F9").Value = 0.641025641025641
$ws.Range("G9").Value = 1.297688566610237
$ws.Range("C10").Value = 39
$ws.Range("D10").Value = 30.73377161304882
$ws.Range("E10").Value = 23.68359261388765
$ws.Range("F10").Value = 3.792409960733627
$ws.Range("G10").Value = 7.677332591670776
$ws.Range("C11").Value = 39
$ws.Range("D11").Value = 24.21550217107311
$ws.Range("E11").Value = 7.378296118301681
$ws.Range("F11").Value = 1.181472935650889
$ws.Range("G11").Value = 2.391766915751601
$ws.Range("C12").Value = 39
$ws.Range("D12").Value = 8.15521326629548
$ws.Range("E12").Value = 7.253684089832208
$ws.Range("F12").Value = 1.161519041590165
$ws.Range("G12").Value = 2.351372369067755
$ws.Range("C13").Value = 39
$ws.Range("D13").Value = 5.982905982905983
$ws.Range("E13").Value = 22.77758577099987
$ws.Range("F13").Value = 3.647332757647231
$ws.Range("G13").Value = 7.383639148426004
$ws.Range("C15").Value = 39
$ws.Range("C16").Value = 39
$ws.Range("D16").Value = 642.2858618830843
$ws.Range("E16").Value = 458.0006948979085
$ws.Range("F16").Value = 73.33880571544907
$ws.Range("G16").Value = 148.4666502786289
$ws.Range("C17").Value = 39
$ws.Range("D17").Value = 1075.61414511012
$ws.Range("E17").Value = 742.217984408955
$ws.Range("F17").Value = 118.8499955643389
$ws.Range("G17").Value = 240.5992374014111
$ws.Range("C18").Value = 39
$ws.Range("D18").Value = 2.128205128205128
$ws.Range("E18").Value = 1.004712244205882
$ws.Range("F18").Value = 0.1608827167700537
$ws.Range("G18").Value = 0.3256900329035991
$ws.Range("C19").Value = 39
$ws.Range("D19").Value = 2.561528466097053
$ws.Range("E19").Value = 0.6080525096128232
$ws.Range("F19").Value = 0.09736632578085137
$ws.Range("G19").Value = 0.197107821672307
$ws.Range("C20").Value = 39
$ws.Range("D20").Value = 688.5745705607274
$ws.Range("E20").Value = 527.8021201861769
$ws.Range("F20").Value = 84.51597907982323
$ws.Range("G20").Value = 171.0936548065002
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 1649.734873685093
$ws.Range("E21").Value = 1777.009504878768
$ws.Range("F21").Value = 284.5492513103291
$ws.Range("G21").Value = 576.0398436981505
$ws.Range("C22").Value = 39
$ws.Range("D22").Value = 60.55790735404162
$ws.Range("E22").Value = 44.14036323945763
$ws.Range("F22").Value = 7.068114873820286
$ws.Range("G22").Value = 14.30865050042117
$ws.Range("C23").Value = 116
$ws.Range("D23").Value = 2.53386116790094
$ws.Range("E23").Value = 3.452382831153333
$ws.Range("F23").Value = 0.320545698673833
$ws.Range("G23").Value = 0.6349393372015498
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 0.5145744069746477
$ws.Range("E24").Value = 0.3493109811986957
$ws.Range("F24").Value = 0.08472035715708043
$ws.Range("G24").Value = 0.1795991340892121
$ws.Range("C25").Value = 116
$ws.Range("D25").Value = 1195.00683632559
$ws.Range("E25").Value = 1614.305504251496
$ws.Range("F25").Value = 149.8845032665289
$ws.Range("G25").Value = 296.892354364954
$ws.Range("C26").Value = 116
$ws.Range("D26").Value = 29.10098522167488
$ws.Range("E26").Value = 31.10144209662779
$ws.Range("F26").Value = 2.887696403963647
$ws.Range("G26").Value = 5.719970813389834
$ws.Range("C27").Value = 116
$ws.Range("D27").Value = 10.43965517241379
$ws.Range("E27").Value = 7.684302229195907
$ws.Range("F27").Value = 0.7134695505526037
$ws.Range("G27").Value = 1.413245866082614
$ws.Range("C28").Value = 116
$ws.Range("D28").Value = 0.4366258174390479
$ws.Range("E28").Value = 5.949938097871621
$ws.Range("F28").Value = 0.5524378836083976
$ws.Range("G28").Value = 1.094273125842998
$ws.Range("C29").Value = 116
$ws.Range("D29").Value = 52.41858237547893
$ws.Range("E29").Value = 35.37845085545433
$ws.Range("F29").Value = 3.284806697891901
$ws.Range("G29").Value = 6.506569878252909
$ws.Range("C30").Value = 116
$ws.Range("D30").Value = 5.373563218390805
$ws.Range("E30").Value = 17.32062829646703
$ws.Range("F30").Value = 1.60817996447573
$ws.Range("G30").Value = 3.185495001085743
$ws.Range("C31").Value = 116
$ws.Range("D31").Value = 20.73291243254585
$ws.Range("E31").Value = 27.77164798173055
$ws.Range("F31").Value = 2.578532781850747
$ws.Range("G31").Value = 5.107577179273602
$ws.Range("C32").Value = 116
$ws.Range("D32").Value = 21.57149253088021
$ws.Range("E32").Value = 7.922046705295242
$ws.Range("F32").Value = 0.7355435709970998
$ws.Range("G32").Value = 1.456970252241554
$ws.Range("C33").Value = 116
$ws.Range("D33").Value = 5.07908320782016
$ws.Range("E33").Value = 7.713824288368586
$ws.Range("F33").Value = 0.7162106049334805
$ws.Range("G33").Value = 1.418675367270831
$ws.Range("C34").Value = 116
$ws.Range("D34").Value = 11.16242474001095
$ws.Range("E34").Value = 23.52615653781688
$ws.Range("F34").Value = 2.184348797148083
$ws.Range("G34").Value = 4.326774569792176
$ws.Range("C36").Value = 116
$ws.Range("D36").Value = 1.944444444444444
$ws.Range("E36").Value = 8.187123152075237
$ws.Range("F36").Value = 0.760155301210891
$ws.Range("G36").Value = 1.505721353048647
$ws.Range("C37").Value = 116
$ws.Range("D37").Value = 468.8803595131669
$ws.Range("E37").Value = 655.8886307950066
$ws.Range("F37").Value = 60.89773055098114
$ws.Range("G37").Value = 120.6266839114974
$ws.Range("C38").Value = 116
$ws.Range("D38").Value = 735.277089864454
$ws.Range("E38").Value = 1090.250733030782
$ws.Range("F38").Value = 101.2272392839649
$ws.Range("G38").Value = 200.5116789388075
$ws.Range("C39").Value = 116
$ws.Range("D39").Value = 2.336206896551724
$ws.Range("E39").Value = 1.509068042475749
$ws.Range("F39").Value = 0.140113450239858
$ws.Range("G39").Value = 0.277537778845198
$ws.Range("C40").Value = 116
$ws.Range("D40").Value = 2.592707400715843
$ws.Range("E40").Value = 0.5471008028821952
$ws.Range("F40").Value = 0.05079703430407292
$ws.Range("G40").Value = 0.1006191486152216
$ws.Range("C41").Value = 116
$ws.Range("D41").Value = 410.0351026777212
$ws.Range("E41").Value = 624.6390854871045
$ws.Range("F41").Value = 57.99628310906614
$ws.Range("G41").Value = 114.8794749384355
$ws.Range("C42").Value = 116
$ws.Range("D42").Value = 1017.715482116955
$ws.Range("E42").Value = 1474.818603196777
$ws.Range("F42").Value = 136.9334696352164
$ws.Range("G42").Value = 271.2388492829598
$ws.Range("C43").Value = 116
$ws.Range("D43").Value = 39.91404339297677
$ws.Range("E43").Value = 50.51568515966026
$ws.Range("F43").Value = 4.690263619484296
$ws.Range("G43").Value = 9.290509547239811
$ws.Range("C44").Value = 71
$ws.Range("D44").Value = 4.26440902099942
$ws.Range("E44").Value = 3.157886917969404
$ws.Range("F44").Value = 0.3747722272896406
$ws.Range("G44").Value = 0.7474596385676053
$ws.Range("C45").Value = 29
$ws.Range("D45").Value = 1.438602542011752
$ws.Range("E45").Value = 0.707232881862319
$ws.Range("F45").Value = 0.1313298491673543
$ws.Range("G45").Value = 0.2690170009653008
$ws.Range("C46").Value = 71
$ws.Range("D46").Value = 1812.428254674726
$ws.Range("E46").Value = 1319.891386286242
$ws.Range("F46").Value = 156.642288805257
$ws.Range("G46").Value = 312.4131940659847
$ws.Range("C47").Value = 71
$ws.Range("D47").Value = 15.4626046879568
$ws.Range("E47").Value = 18.18468749248531
$ws.Range("F47").Value = 2.158125357606874
$ws.Range("G47").Value = 4.304245305065613
$ws.Range("C48").Value = 71
$ws.Range("D48").Value = 15.47183098591549
$ws.Range("E48").Value = 10.2946659757361
$ws.Range("F48").Value = 1.221752075723593
$ws.Range("G48").Value = 2.436707681206614
$ws.Range("C49").Value = 71
$ws.Range("D49").Value = 3.272684639357442
$ws.Range("E49").Value = 4.262454029584707
$ws.Range("F49").Value = 0.505860226120561
$ws.Range("G49").Value = 1.008906408343811
$ws.Range("C50").Value = 71
$ws.Range("D50").Value = 72.67915183408142
$ws.Range("E50").Value = 24.07693339319469
$ws.Range("F50").Value = 2.857406293659856
$ws.Range("G50").Value = 5.698917155483772
$ws.Range("C51").Value = 71
$ws.Range("C52").Value = 71
$ws.Range("D52").Value = 37.16285098036012
$ws.Range("E52").Value = 26.66069555970455
$ws.Range("F52").Value = 3.164042448494766
$ws.Range("G52").Value = 6.310483682497334
$ws.Range("C53").Value = 71
$ws.Range("D53").Value = 27.03648127617942
$ws.Range("E53").Value = 4.553386938601863
$ws.Range("F53").Value = 0.5403876101392152
$ws.Range("G53").Value = 1.07776910440299
$ws.Range("C54").Value = 71
$ws.Range("D54").Value = 10.01901791274749
$ws.Range("E54").Value = 7.745578461584078
$ws.Range("F54").Value = 0.9192310450309347
$ws.Range("G54").Value = 1.833348510501907
$ws.Range("C55").Value = 71
$ws.Range("D55").Value = 10.71806909835079
$ws.Range("E55").Value = 17.47747763799126
$ws.Range("F55").Value = 2.074194989226126
$ws.Range("G55").Value = 4.136851463562421
$ws.Range("C57").Value = 71
$ws.Range("D57").Value = 1.140174379610999
$ws.Range("E57").Value = 5.774332396965119
$ws.Range("F57").Value = 0.6852871777032751
$ws.Range("G57").Value = 1.366762179432347
$ws.Range("C58").Value = 71
$ws.Range("D58").Value = 781.0129051957277
$ws.Range("E58").Value = 581.1968136737165
$ws.Range("F58").Value = 68.97537182686773
$ws.Range("G58").Value = 137.5670413697217
$ws.Range("C59").Value = 71
$ws.Range("D59").Value = 1117.877447111412
$ws.Range("E59").Value = 881.2847071199818
$ws.Range("F59").Value = 104.5892526056736
$ws.Range("G59").Value = 208.5966868891667
$ws.Range("C60").Value = 71
$ws.Range("D60").Value = 2.619718309859155
$ws.Range("E60").Value = 1.561740761212112
$ws.Range("F60").Value = 0.1853445290260241
$ws.Range("G60").Value = 0.3696580071532544
$ws.Range("C61").Value = 71
$ws.Range("D61").Value = 2.288079340545889
$ws.Range("E61").Value = 0.3102880286550747
$ws.Range("F61").Value = 0.03682441411649715
$ws.Range("G61").Value = 0.07344397813317266
$ws.Range("C62").Value = 71
$ws.Range("D62").Value = 769.1558775632627
$ws.Range("E62").Value = 588.2224022443697
$ws.Range("F62").Value = 69.80915579223449
$ws.Range("G62").Value = 139.2299710534489
$ws.Range("C63").Value = 71
$ws.Range("D63").Value = 1458.726533716801
$ws.Range("E63").Value = 1297.318498732159
$ws.Range("F63").Value = 153.9633799130912
$ws.Range("G63").Value = 307.0702787523956
$ws.Range("C64").Value = 71
$ws.Range("D64").Value = 71.96579674413367
$ws.Range("E64").Value = 51.02329362506878
$ws.Range("F64").Value = 6.055350901487099
$ws.Range("G64").Value = 12.07701656272298

# Row 56 ("SM + Traps" / "Percent of catch = mature") now has N = 0 and no
# stats (underlying sample group became empty), so D56:G56 are cleared
# entirely (not just zeroed) and C56 is set to 0.
$ws.Range("D56:G56").ClearContents()
$ws.Range("C56").Value = 0
